# Data refresh for the Goblin_Profits workbook (per-sheet market-price / profit
# columns H:N). Mirrors a scheduled-runner update: some rows also gain or lose
# cells entirely (Excel drops a numeric cell whose value becomes blank), which
# is reproduced below by assigning "" to those specific cells.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 82.166664
$ws.Range("I4").Value = 82.166664
$ws.Range("K4").Value = 82.166664
$ws.Range("M4").Value = 31.833336

$ws.Range("H11").Value = 194.91667
$ws.Range("I11").Value = 194.91667
$ws.Range("K11").Value = 194.91667
$ws.Range("M11").Value = -54.91667000000001

$ws.Range("H32").Value = 3000
$ws.Range("I32").Value = 1000
$ws.Range("J32").Value = 5000
$ws.Range("K32").Value = 1000
$ws.Range("L32").Value = 5000
$ws.Range("M32").Value = -674
$ws.Range("N32").Value = -5652

$ws.Range("H40").Value = 1169.8
$ws.Range("I40").Value = 1237.5
$ws.Range("K40").Value = 1237.5
$ws.Range("M40").Value = -1062.5

$ws.Range("H57").Value = 30380.953
$ws.Range("J57").Value = 30380.953
$ws.Range("L57").Value = 91142.859
$ws.Range("N57").Value = -92140.859

$ws.Range("H62").Value = 58740.668
$ws.Range("I62").Value = 91145.22
$ws.Range("J62").Value = 10133.833
$ws.Range("K62").Value = 91145.22
$ws.Range("L62").Value = 10133.833
$ws.Range("M62").Value = -90521.22
$ws.Range("N62").Value = -11381.833

$ws.Range("H65").Value = 58740.668
$ws.Range("I65").Value = 91145.22
$ws.Range("J65").Value = 10133.833
$ws.Range("K65").Value = 455726.1
$ws.Range("L65").Value = 50669.165
$ws.Range("M65").Value = -452606.1
$ws.Range("N65").Value = -56909.165

$ws.Range("H107").Value = 748.1667
$ws.Range("I107").Value = 970.25
$ws.Range("J107").Value = 304
$ws.Range("K107").Value = 970.25
$ws.Range("L107").Value = 304
$ws.Range("M107").Value = 949.75
$ws.Range("N107").Value = -4144

$ws.Range("H131").Value = 9208.200000000001
$ws.Range("J131").Value = 9283.857
$ws.Range("L131").Value = 27851.571
$ws.Range("N131").Value = -37931.571

$ws.Range("H132").Value = 1784.4333
$ws.Range("I132").Value = 1538.3043
$ws.Range("J132").Value = 2593.1428
$ws.Range("K132").Value = 4614.9129
$ws.Range("L132").Value = 7779.428400000001
$ws.Range("M132").Value = -2084.9129
$ws.Range("N132").Value = -12839.4284

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2500.5454
$ws.Range("I74").Value = 2738.625
$ws.Range("J74").Value = 1865.6666
$ws.Range("K74").Value = 2738.625
$ws.Range("L74").Value = 1865.6666
$ws.Range("M74").Value = -1864.625
$ws.Range("N74").Value = -3613.6666

$ws.Range("H77").Value = 2500.5454
$ws.Range("I77").Value = 2738.625
$ws.Range("J77").Value = 1865.6666
$ws.Range("K77").Value = 13693.125
$ws.Range("L77").Value = 9328.333000000001
$ws.Range("M77").Value = -9325.125
$ws.Range("N77").Value = -18064.333

$ws.Range("H97").Value = 1534.6666
$ws.Range("I97").Value = 1508.8572
$ws.Range("J97").Value = 1625
$ws.Range("K97").Value = 1508.8572
$ws.Range("L97").Value = 1625
$ws.Range("M97").Value = -1012.8572
$ws.Range("N97").Value = -2617

$ws.Range("H122").Value = 2868.8333
$ws.Range("I122").Value = 2749
$ws.Range("J122").Value = 2892.8
$ws.Range("K122").Value = 8247
$ws.Range("L122").Value = 8678.400000000001
$ws.Range("M122").Value = -5797
$ws.Range("N122").Value = -13578.4

$ws.Range("H132").Value = 1956.12
$ws.Range("I132").Value = 1937.619
$ws.Range("J132").Value = 2053.25
$ws.Range("K132").Value = 5812.857
$ws.Range("L132").Value = 6159.75
$ws.Range("M132").Value = -3282.857
$ws.Range("N132").Value = -11219.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 3312.44
$ws.Range("I94").Value = 3120.6316
$ws.Range("J94").Value = 3919.8333
$ws.Range("K94").Value = 3120.6316
$ws.Range("L94").Value = 3919.8333
$ws.Range("M94").Value = -2669.6316
$ws.Range("N94").Value = -4821.8333

$ws.Range("H96").Value = 26092.625
$ws.Range("I96").Value = 7749.2
$ws.Range("K96").Value = 7749.2
$ws.Range("M96").Value = -5003.2

$ws.Range("H134").Value = 4345.273
$ws.Range("I134").Value = 4474.75
$ws.Range("K134").Value = 13424.25
$ws.Range("M134").Value = -10889.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 5051.909
$ws.Range("I132").Value = 5357.2
$ws.Range("J132").Value = 1999
$ws.Range("K132").Value = 16071.6
$ws.Range("L132").Value = 5997
$ws.Range("M132").Value = -13541.6
$ws.Range("N132").Value = -11057

$ws.Range("H134").Value = 3917.9375
$ws.Range("I134").Value = 3712.4666
$ws.Range("J134").Value = 7000
$ws.Range("K134").Value = 11137.3998
$ws.Range("L134").Value = 21000
$ws.Range("M134").Value = -8602.399800000001
$ws.Range("N134").Value = -26070

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 6734291.5
$ws.Range("I17").Value = 16666945
$ws.Range("K17").Value = 50000835
$ws.Range("M17").Value = -50000666

$ws.Range("H121").Value = 866.4286
$ws.Range("I121").Value = 718.75
$ws.Range("J121").Value = 925.5
$ws.Range("K121").Value = 2156.25
$ws.Range("L121").Value = 2776.5
$ws.Range("M121").Value = -846.25
$ws.Range("N121").Value = -5396.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 2749.5
$ws.Range("I5").Value = 1999
$ws.Range("K5").Value = 1999
$ws.Range("M5").Value = -1887

$ws.Range("H55").Value = 34999.25
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 34999.25
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = ""
$ws.Range("M55").Value = 34999.25
$ws.Range("N55").Value = -35653.25

$ws.Range("H70").Value = 1000
$ws.Range("I70").Value = 1000
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 1000
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = ""
$ws.Range("N70").Value = -730

$ws.Range("H73").Value = 1000
$ws.Range("I73").Value = 1000
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 1000
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = ""
$ws.Range("N73").Value = -64

$ws.Range("H126").Value = 2596.9473
$ws.Range("I126").Value = 2609.5625
$ws.Range("J126").Value = 2529.6667
$ws.Range("K126").Value = 7828.6875
$ws.Range("L126").Value = 7589.000100000001
$ws.Range("M126").Value = -5358.6875
$ws.Range("N126").Value = -12529.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2466.9333
$ws.Range("I132").Value = 2263.7036
$ws.Range("J132").Value = 4296
$ws.Range("K132").Value = 6791.110799999999
$ws.Range("L132").Value = 12888
$ws.Range("M132").Value = -4261.110799999999
$ws.Range("N132").Value = -17948

$ws.Range("H133").Value = 0
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = ""
$ws.Range("M133").Value = ""
$ws.Range("N133").Value = 0

$ws.Range("H136").Value = 6160.4
$ws.Range("I136").Value = 6199.25
$ws.Range("J136").Value = 6005
$ws.Range("K136").Value = 18597.75
$ws.Range("L136").Value = 18015
$ws.Range("M136").Value = -16047.75
$ws.Range("N136").Value = -23115

$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = ""
$ws.Range("N139").Value = 0

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 10000
$ws.Range("I8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("M8").Value = ""

$ws.Range("H122").Value = 3646.5789
$ws.Range("I122").Value = 3890.8235
$ws.Range("J122").Value = 1570.5
$ws.Range("K122").Value = 11672.4705
$ws.Range("L122").Value = 4711.5
$ws.Range("M122").Value = -9222.470499999999
$ws.Range("N122").Value = -9611.5

$ws.Range("H132").Value = 1381.8077
$ws.Range("I132").Value = 1267.08
$ws.Range("K132").Value = 3801.24
$ws.Range("M132").Value = -1271.24
